# Consolidate adjacent text-run nodes that share identical formatting into
# a single <a:r> per "word", matching the OOXML produced by the updated
# PowerPoint writer (moves the trailing space into the preceding run's
# text instead of keeping a standalone space-only run).
#
# Runs() rounds to already-merged COM run boundaries, so we use
# Characters(start, length) to address the exact character ranges that
# back each individual <a:r> in the underlying OOXML, then delete the
# now-redundant standalone space runs (setting Text = "" removes a run).

$p = $ppt.ActivePresentation

# --- Slide 1 title: "Header" + " " + "with" + " " + "inline code"(Courier)
#     -> "Header " + "with " + "inline code"(Courier)
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange

$tr1.Characters(1, 6).Text = "Header "   # "Header" -> "Header "
$tr1.Characters(8, 1).Text = ""          # drop the now-redundant " " run

$tr1.Characters(8, 4).Text = "with "     # "with" -> "with "
$tr1.Characters(13, 1).Text = ""         # drop the now-redundant " " run

# --- Slide 2 title: "Syntax" + " " + "highlighting"
#     -> "Syntax " + "highlighting"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange

$tr2.Characters(1, 6).Text = "Syntax "   # "Syntax" -> "Syntax "
$tr2.Characters(8, 1).Text = ""          # drop the now-redundant " " run

# --- Slide 3 title: "Two" + " " + "column" + " " + "slide"
#     -> "Two " + "column " + "slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange

$tr3.Characters(1, 3).Text = "Two "      # "Two" -> "Two "
$tr3.Characters(5, 1).Text = ""          # drop the now-redundant " " run

$tr3.Characters(5, 6).Text = "column "   # "column" -> "column "
$tr3.Characters(12, 1).Text = ""         # drop the now-redundant " " run
